$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.892.89"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "1.624.25"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.66%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.01"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0604"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0882"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("D12").Value = "1.855.47"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.639.13"
$ws.Range("E13").Value = "  +0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.97"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.49"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.67%  "
$ws.Range("D17").Value = "27.874.90"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0715"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.93"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.07"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D33").Value = "1.421.10"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.971"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  -0.80%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  -2.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.02"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("D46").Value = "1.764.78"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.25"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.24%  "
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1000"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.34%  "
$ws.Range("E51").Value = "  -0.34%  "

Write-Output "Applied 81 changes"